$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with merged content
$ws.Range("A2").Value = '[{"role": "assistant", "content": "Chào cậu! Bữa nay thầy giáo robot trên sao HỎA dạy tớ cách luyện nói theo nhịp như… hát rap đó! Mỗi cụm từ giống như một đoạn nhạc nhỏ, có nhịp điệu riêng."}, {"role": "assistant", "content": "Cứ theo nhịp là nói được tiếng Anh hay hơn liền! Giờ tụi mình thử cùng luyện nha. Cậu đã sẵn sàng chưa?"}, {"role": "user", "content": "[NEXT]"}]'
$ws.Range("C2").Value = "Chào cậu! Bữa nay thầy giáo robot trên sao HỎA dạy tớ cách luyện nói theo nhịp như… hát rap đó! Mỗi cụm từ giống như một đoạn nhạc nhỏ, có nhịp điệu riêng."
$ws.Range("D2").Value = "Đã rồi! Cậu chuẩn bị sẵn sàng rồi!"

# Delete rows 3 through 7 (shift remaining rows up / remove entirely)
$ws.Range("A3:D7").EntireRow.Delete()
